$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.471.95"
$ws.Range("E2").Value = "  -0.58%  "
$ws.Range("D3").Value = "1.836.68"
$ws.Range("E3").Value = "  -0.83%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "260.79"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.51%  "
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5370"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +2.29%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3012"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -7.37%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06864"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.89%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "17.56"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -7.58%  "
$ws.Range("B11").Value = "WrappedEther"
$ws.Range("C11").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D11").Value = "1.853.85"
$ws.Range("E11").Value = "  +0.20%  "
$ws.Range("B12").Value = "Polygon"
$ws.Range("C12").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7367"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -5.76%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07147"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -8.39%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "89.47"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.91%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.977"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.02%  "
$ws.Range("E16").Value = "  +0.02%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.80"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -1.57%  "
$ws.Range("E18").Value = "  +0.08%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007870"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.42%  "
$ws.Range("D20").Value = "26.496.38"
$ws.Range("E20").Value = "  -0.55%  "
$ws.Range("D21").Value = "2.076.94"
$ws.Range("E21").Value = "  -0.25%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.581"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.35%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.959"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.90%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.216"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -3.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "142.78"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.07%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.182"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.688"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.27%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.92"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.89%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "110.63"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -1.32%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.213"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.22%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08817"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.96%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.012"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -2.63%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04804"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.80%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.924"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +1.42%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7265"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.41%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.129"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.19%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.094"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.71%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.264"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.21%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01704"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -5.13%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.4709"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -3.54%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9041"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "107.59"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -3.07%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.872"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.96%  "
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.366"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -4.11%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.997"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.78%  "
$ws.Range("B47").Value = "Algorand"
$ws.Range("C47").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1236"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.05%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4049"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -3.93%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "34.69"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.19%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05775"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.96%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.8889"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.07%  "
